# Updates cryptos list values to match the latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '66.563.03'
$ws.Range("E2").Value = '  -3.55%  '

# Row 3
$ws.Range("D3").Value = '3.559.25'
$ws.Range("E3").Value = '  -4.28%  '

# Row 5
$ws.Range("D5").Value = '''571.35'
$ws.Range("E5").Value = '  -6.77%  '

# Row 6
$ws.Range("D6").Value = '''187.24'
$ws.Range("E6").Value = '  -2.36%  '

# Row 7
$ws.Range("D7").Value = '3.554.85'
$ws.Range("E7").Value = '  -4.25%  '

# Row 8
$ws.Range("D8").Value = '''0.610'
$ws.Range("E8").Value = '  -4.61%  '

# Row 9
$ws.Range("D9").Value = '''0.998'
$ws.Range("E9").Value = '  -0.24%  '

# Row 10
$ws.Range("D10").Value = '''0.666'
$ws.Range("E10").Value = '  -8.48%  '

# Row 11
$ws.Range("D11").Value = '''0.146'
$ws.Range("E11").Value = '  -9.72%  '

# Row 12
$ws.Range("D12").Value = '''55.05'
$ws.Range("E12").Value = '  -9.12%  '

# Row 13
$ws.Range("D13").Value = '''0.0000259'
$ws.Range("E13").Value = '  -11.15%  '

# Row 14
$ws.Range("D14").Value = '''9.75'
$ws.Range("E14").Value = '  -7.39%  '

# Row 15
$ws.Range("D15").Value = '4.116.41'
$ws.Range("E15").Value = '  -4.73%  '

# Row 16
$ws.Range("D16").Value = '3.544.87'
$ws.Range("E16").Value = '  -4.96%  '

# Row 17
$ws.Range("E17").Value = '  -1.44%  '

# Row 18
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").Value = '''18.17'
$ws.Range("E18").Value = '  -6.94%  '

# Row 19
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '66.477.33'
$ws.Range("E19").Value = '  -3.55%  '

# Row 20
$ws.Range("D20").Value = '''12.04'
$ws.Range("E20").Value = '  -7.31%  '

# Row 21
$ws.Range("D21").Value = '''1.06'
$ws.Range("E21").Value = '  -8.90%  '

# Row 22
$ws.Range("D22").Value = '''387.71'
$ws.Range("E22").Value = '  -6.20%  '

# Row 23
$ws.Range("D23").Value = '''4.19'
$ws.Range("E23").Value = '  -8.53%  '

# Row 24
$ws.Range("D24").Value = '''85.18'
$ws.Range("E24").Value = '  -5.63%  '

# Row 25
$ws.Range("D25").Value = '''11.09'
$ws.Range("E25").Value = '  -3.32%  '

# Row 26
$ws.Range("D26").Value = '''2.91'
$ws.Range("E26").Value = '  -5.99%  '

# Row 27
$ws.Range("D27").Value = '''12.33'
$ws.Range("E27").Value = '  -5.33%  '

# Row 28
$ws.Range("D28").Value = '''6.05'
$ws.Range("E28").Value = '  +0.38%  '

# Row 29
$ws.Range("D29").Value = '''3.53'
$ws.Range("E29").Value = '  -7.73%  '

# Row 30
$ws.Range("D30").Value = '''8.85'
$ws.Range("E30").Value = '  -9.82%  '

# Row 31
$ws.Range("D31").Value = '''7.56'
$ws.Range("E31").Value = '  -4.03%  '

# Row 32
$ws.Range("D32").Value = '''30.87'
$ws.Range("E32").Value = '  -6.41%  '

# Row 33
$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").Value = '''624.00'
$ws.Range("E33").Value = '  -4.03%  '

# Row 34
$ws.Range("B34").Value = 'Cosmos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D34").Value = '''12.15'
$ws.Range("E34").Value = '  -5.05%  '

# Row 35
$ws.Range("D35").Value = '''0.114'
$ws.Range("E35").Value = '  -7.50%  '

# Row 36
$ws.Range("D36").Value = '''63.22'
$ws.Range("E36").Value = '  -5.48%  '

# Row 37
$ws.Range("D37").Value = '''41.60'
$ws.Range("E37").Value = '  -10.17%  '

# Row 38
$ws.Range("D38").Value = '''0.404'
$ws.Range("E38").Value = '  -3.03%  '

# Row 39
$ws.Range("E39").Value = '  +0.12%  '

# Row 40
$ws.Range("D40").Value = '0.0₃0736'
$ws.Range("E40").Value = '  -12.27%  '

# Row 41
$ws.Range("D41").Value = '''0.133'
$ws.Range("E41").Value = '  -5.84%  '

# Row 42
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '3.105.64'
$ws.Range("E42").Value = '  +6.90%  '

# Row 43
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").Value = '''0.997'
$ws.Range("E43").Value = '  -0.34%  '

# Row 44
$ws.Range("E44").Value = '  -6.38%  '

# Row 45
$ws.Range("D45").Value = '''2.54'
$ws.Range("E45").Value = '  -3.75%  '

# Row 46
$ws.Range("D46").Value = '''0.0407'
$ws.Range("E46").Value = '  -9.50%  '

# Row 47
$ws.Range("D47").Value = '''3.09'
$ws.Range("E47").Value = '  +0.44%  '

# Row 48
$ws.Range("E48").Value = '  -7.93%  '

# Row 49
$ws.Range("D49").Value = '''138.26'
$ws.Range("E49").Value = '  -3.61%  '

# Row 50
$ws.Range("D50").Value = '''8.41'
$ws.Range("E50").Value = '  -9.34%  '

# Row 51
$ws.Range("D51").Value = '''2.74'
$ws.Range("E51").Value = '  -1.66%  '
